# Regression Test Plan - iOS PROD-TEST sheet:
# Remove the "services duty cycle" warning from the "Time and Date Testing
# Was Complete" instructions cell and replace it with a short, plain
# statement that testing can happen at any time/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A3")

# Retype the cell with the simplified wording (this is what previously was
# four differently-formatted text runs - bold lead-in, bold+underlined time
# window, bold closing paren, plain duty-cycle pointer - collapsing them
# into a single, unformatted run).
$cell.Value = "Time and Date Testing Was Complete (Testing can be completed any time/date):"

# The old text was bold; the new, shorter instruction is typed in plain
# (non-bold) style.
$cell.Font.Bold = $False

# The cell no longer needs four wrapped lines, so the row shrinks back down
# to (roughly) a single line of text.
$ws.Range("A3:U3").RowHeight = 15.75

# Leave the selection sitting on the edited cell.
$cell.Select() | Out-Null
